# Updates the "Estado de Cuenta" workbook with the refreshed worker/debt
# database described in the commit ("Actualiza base de datos EC y agrega
# parte 1 de nuevos estado de cuenta").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header summary values -------------------------------------------------
$ws.Range("E11").Value = 1092000
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 12

# --- Remove the two extra rows that belonged to the old last worker -------
# (row 38 keeps the "last row" border styling and becomes the new row 36)
$ws.Rows("36:37").Delete() | Out-Null

# --- Row 16: GLAUCO SARKAR HERNANDEZ ---------------------------------------
$ws.Range("C16").Value = "73581395"
$ws.Range("D16").Value = "GLAUCO SARKAR HERNANDEZ"
$ws.Range("E16").Value = "2408"
$ws.Range("F16").Value = 52000
$ws.Range("G16").Value = 1300000

# --- Row 17: CARLOS ALBERTO PABA ALMANZA -----------------------------------
$ws.Range("C17").Value = "1007574877"
$ws.Range("D17").Value = "CARLOS ALBERTO PABA ALMANZA"
$ws.Range("E17").Value = "2410"
$ws.Range("F17").Value = 52000
$ws.Range("G17").Value = 1300000

# --- Row 18: CARLOS ALBERTO PABA ALMANZA -----------------------------------
$ws.Range("C18").Value = "1007574877"
$ws.Range("D18").Value = "CARLOS ALBERTO PABA ALMANZA"
$ws.Range("E18").Value = "2411"
$ws.Range("F18").Value = 52000
$ws.Range("G18").Value = 1300000

# --- Row 19: DANUBIS HERNANDEZ CABRERA --------------------------------------
$ws.Range("C19").Value = "1047437485"
$ws.Range("D19").Value = "DANUBIS HERNANDEZ CABRERA"
$ws.Range("E19").Value = "2412"
$ws.Range("F19").Value = 52000
$ws.Range("G19").Value = 1300000

# --- Row 20: CARLOS ALBERTO PABA ALMANZA -----------------------------------
$ws.Range("C20").Value = "1007574877"
$ws.Range("D20").Value = "CARLOS ALBERTO PABA ALMANZA"
$ws.Range("E20").Value = "2412"
$ws.Range("F20").Value = 52000
$ws.Range("G20").Value = 1300000

# --- Row 21: DANUBIS HERNANDEZ CABRERA --------------------------------------
$ws.Range("C21").Value = "1047437485"
$ws.Range("D21").Value = "DANUBIS HERNANDEZ CABRERA"
$ws.Range("E21").Value = "2501"
$ws.Range("F21").Value = 52000
$ws.Range("G21").Value = 1300000

# --- Row 22: CARLOS ALBERTO PABA ALMANZA -----------------------------------
$ws.Range("C22").Value = "1007574877"
$ws.Range("D22").Value = "CARLOS ALBERTO PABA ALMANZA"
$ws.Range("E22").Value = "2501"
$ws.Range("F22").Value = 52000
$ws.Range("G22").Value = 1300000

# --- Row 23: DANUBIS HERNANDEZ CABRERA --------------------------------------
$ws.Range("C23").Value = "1047437485"
$ws.Range("D23").Value = "DANUBIS HERNANDEZ CABRERA"
$ws.Range("E23").Value = "2502"
$ws.Range("F23").Value = 52000
$ws.Range("G23").Value = 1300000

# --- Row 24: CARLOS ALBERTO PABA ALMANZA -----------------------------------
$ws.Range("C24").Value = "1007574877"
$ws.Range("D24").Value = "CARLOS ALBERTO PABA ALMANZA"
$ws.Range("E24").Value = "2502"
$ws.Range("F24").Value = 52000
$ws.Range("G24").Value = 1300000

# --- Row 25: DANUBIS HERNANDEZ CABRERA --------------------------------------
$ws.Range("C25").Value = "1047437485"
$ws.Range("D25").Value = "DANUBIS HERNANDEZ CABRERA"
$ws.Range("E25").Value = "2503"
$ws.Range("F25").Value = 52000
$ws.Range("G25").Value = 1300000

# --- Row 26: CARLOS ALBERTO PABA ALMANZA -----------------------------------
$ws.Range("C26").Value = "1007574877"
$ws.Range("D26").Value = "CARLOS ALBERTO PABA ALMANZA"
$ws.Range("E26").Value = "2503"
$ws.Range("F26").Value = 52000
$ws.Range("G26").Value = 1300000

# --- Row 27: DANUBIS HERNANDEZ CABRERA --------------------------------------
$ws.Range("C27").Value = "1047437485"
$ws.Range("D27").Value = "DANUBIS HERNANDEZ CABRERA"
$ws.Range("E27").Value = "2504"
$ws.Range("F27").Value = 52000
$ws.Range("G27").Value = 1300000

# --- Row 28: CARLOS ALBERTO PABA ALMANZA -----------------------------------
$ws.Range("C28").Value = "1007574877"
$ws.Range("D28").Value = "CARLOS ALBERTO PABA ALMANZA"
$ws.Range("E28").Value = "2504"
$ws.Range("F28").Value = 52000
$ws.Range("G28").Value = 1300000

# --- Row 29: DANUBIS HERNANDEZ CABRERA --------------------------------------
$ws.Range("C29").Value = "1047437485"
$ws.Range("D29").Value = "DANUBIS HERNANDEZ CABRERA"
$ws.Range("E29").Value = "2505"
$ws.Range("F29").Value = 52000
$ws.Range("G29").Value = 1300000

# --- Row 30: CARLOS ALBERTO PABA ALMANZA -----------------------------------
$ws.Range("C30").Value = "1007574877"
$ws.Range("D30").Value = "CARLOS ALBERTO PABA ALMANZA"
$ws.Range("E30").Value = "2505"
$ws.Range("F30").Value = 52000
$ws.Range("G30").Value = 1300000

# --- Row 31: DANUBIS HERNANDEZ CABRERA --------------------------------------
$ws.Range("C31").Value = "1047437485"
$ws.Range("D31").Value = "DANUBIS HERNANDEZ CABRERA"
$ws.Range("E31").Value = "2506"
$ws.Range("F31").Value = 52000
$ws.Range("G31").Value = 1300000

# --- Row 32: CARLOS ALBERTO PABA ALMANZA -----------------------------------
$ws.Range("C32").Value = "1007574877"
$ws.Range("D32").Value = "CARLOS ALBERTO PABA ALMANZA"
$ws.Range("E32").Value = "2506"
$ws.Range("F32").Value = 52000
$ws.Range("G32").Value = 1300000

# --- Row 33: DANUBIS HERNANDEZ CABRERA --------------------------------------
$ws.Range("C33").Value = "1047437485"
$ws.Range("D33").Value = "DANUBIS HERNANDEZ CABRERA"
$ws.Range("E33").Value = "2507"
$ws.Range("F33").Value = 52000
$ws.Range("G33").Value = 1300000

# --- Row 34: CARLOS ALBERTO PABA ALMANZA -----------------------------------
$ws.Range("C34").Value = "1007574877"
$ws.Range("D34").Value = "CARLOS ALBERTO PABA ALMANZA"
$ws.Range("E34").Value = "2507"
$ws.Range("F34").Value = 52000
$ws.Range("G34").Value = 1300000

# --- Row 35: DANUBIS HERNANDEZ CABRERA --------------------------------------
$ws.Range("C35").Value = "1047437485"
$ws.Range("D35").Value = "DANUBIS HERNANDEZ CABRERA"
$ws.Range("E35").Value = "2508"
$ws.Range("F35").Value = 52000
$ws.Range("G35").Value = 1300000

# --- Row 36 (was row 38 before the delete above): CARLOS ALBERTO PABA ALMANZA
$ws.Range("C36").Value = "1007574877"
$ws.Range("D36").Value = "CARLOS ALBERTO PABA ALMANZA"
$ws.Range("E36").Value = "2508"
$ws.Range("F36").Value = 52000
$ws.Range("G36").Value = 1300000

# --- Column D width shrinks now that the longest name was removed ---------
# (Excel recalculates the "best fit" width for the shorter names that
# remain; set explicitly to match Excel's own best-fit metric.)
$ws.Columns("D").ColumnWidth = 31.6328125
